$p = $ppt.ActivePresentation

$newDate = "3/8/2022"

function Update-DatePlaceholder($shapes) {
    for ($i = 1; $i -le $shapes.Count; $i++) {
        $sh = $shapes.Item($i)
        if ($sh.HasTextFrame) {
            $ph = $sh.PlaceholderFormat
            if ($ph -ne $null -and $ph.Type -eq 16) {
                $sh.TextFrame.TextRange.Text = $newDate
            }
        }
    }
}

# Update the slide master's Date placeholder.
$master = $p.SlideMaster
Update-DatePlaceholder $master.Shapes

# Update the Date placeholder on every slide layout (custom layout).
for ($L = 1; $L -le $master.CustomLayouts.Count; $L++) {
    $cl = $master.CustomLayouts.Item($L)
    Update-DatePlaceholder $cl.Shapes
}
